$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# The "course" column (C) was mis-entered as "slottsskogen" for every result
# row. Fix it to the correct course name, "Teleborgs discgolfbana", for all
# data rows (2-46).
for ($r = 2; $r -le 46; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq "slottsskogen") {
        $cell.Value = "Teleborgs discgolfbana"
    }
}

# Update the view/selection state left over from the editing session: the
# window had scrolled back up near the top of the sheet, with a different
# cell selected.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J17").Select()
